# Commit: "delete duplicate 2  with linkedlist"
# Adds a new row (20th problem entry) to the "链表" (linked list) worksheet,
# describing LeetCode #82 (Remove Duplicates from Sorted List II).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("链表")

# --- widen columns C (题目) and D (解题方法) to fit the new, longer content ---
$ws.Columns.Item(3).ColumnWidth = 43.714285714285715
$ws.Columns.Item(4).ColumnWidth = 75.14285714285714

# --- new row content (row 21) ---
$title = @'
给定一个排序链表，删除所有含有重复数字的节点，只保留原始链表中 没有重复出现 的数字。 
 示例 1: 
 输入: 1->2->3->3->4->4->5
输出: 1->2->5
 示例 2: 
 输入: 1->1->1->2->3
输出: 2->3 
'@

$method = @'
1 创建solder节点，便于处理头节点重复的情况
2 cur指针从solder开始，
3 cur的next节点与next.next节点的值对比
     如果不相等，就迭代cur=cur.next
     如果相等，此时cur在相等节点之前，tmp指针作为重复节点的第一个元素迭代
     tmp节点与tmp.next节点是否相等
            如果相等，就继续移动tmp
            如果不想等，cur指向tmp的下一个节点【这个节点就是重复节点链后的第一个节点】
     循环条件是tmp的值是否与tmp.next的值相等
4 循环结束条件是cur的next节点与next.next节点不为空
'@

$keywords = @'
双指针
快慢指针

'@

$ws.Range("A21").Value = 20
$ws.Range("B21").Value = 82
$ws.Range("C21").Value = $title
$ws.Range("D21").Value = $method
$ws.Range("E21").Value = $keywords

$ws.Rows.Item(21).RowHeight = 240

# --- move the view / selection the way the author left it ---
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F26").Select()
